$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D, shifting existing D:K data to E:L.
$ws.Columns("D:D").Insert()

# Copy number formats/styles from the (now shifted) old column, currently E, into the
# newly inserted blank column D so the new column matches the original column's styling.
# Restrict to the sheet's actual D:K data blocks (done as three separate single-area
# operations, since this engine treats comma-joined multi-area ranges as one bounding
# box for Copy/PasteSpecial) so we don't materialize stray formatted cells in label-only
# rows that never had D:K data (rows 5, 6, 37, 79).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the newest reporting period's values.
$values = @{
    7  = 43465
    8  = 972300
    9  = 704700
    10 = 267700
    12 = "NA"
    13 = 0
    14 = 4100
    15 = 0
    17 = 1006400
    18 = -34100
    20 = 0
    21 = -18500
    22 = "NA"
    23 = -34100
    24 = 11800
    25 = 0
    26 = -45900
    27 = -45900
    28 = 0
    29 = 2000
    30 = 0
    31 = 0
    32 = 0
    33 = -43900
    34 = 0
    35 = -43900
    38 = 43465
    41 = 116900
    42 = 0
    43 = 138500
    44 = 128100
    45 = 13900
    46 = 397400
    47 = 2300
    48 = 74300
    49 = 407800
    50 = 0
    51 = 0
    52 = 4100
    53 = 0
    54 = 885900
    57 = 92500
    58 = 2100
    59 = 103600
    60 = 198200
    61 = 253500
    62 = 75000
    63 = 0
    64 = 0
    65 = 0
    66 = 526700
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 142400
    73 = 0
    74 = 0
    75 = 0
    76 = 359100
    77 = 0
    80 = 43465
    81 = -43900
    83 = 15600
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = -46400
    91 = -9800
    92 = 0
    93 = 0
    94 = -6400
    96 = -1600
    97 = 0
    98 = 0
    99 = 0
    100 = -2900
    101 = -3900
    102 = -59600
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}
